# Edit the document per the target diff:
#  1. Append two trailing spaces to the first paragraph's text, then add a
#     red-colored annotation "(This is a change – Version for main branch)"
#     split across three runs (mirrors the original authoring runs).
#  2. Append a new, otherwise-empty paragraph at the very end of the body
#     (before the section break) shaded with fill color F9F9F9.

$d = $word.ActiveDocument

# --- 1. First paragraph: trailing spaces + red annotation -------------------
$p1 = $d.Paragraphs(1).Range

$insertionPoint = $p1.End - 1
$spaceRange = $d.Range($insertionPoint, $insertionPoint)
$spaceRange.InsertAfter("  ")

$redParts = @(
  [string][char]0x0028 + "This is a change " + [string][char]0x2013 + " Ve",
  "rsion for main branch",
  [string][char]0x0029
)

foreach ($part in $redParts) {
  $ip = $p1.End - 1
  $insertRange = $d.Range($ip, $ip)
  $insertRange.InsertAfter($part)

  $runStart = $ip
  $runEnd = $p1.End - 1
  $newRunRange = $d.Range($runStart, $runEnd)
  $newRunRange.Font.Color = 255
}

# --- 2. New shaded empty paragraph at the end of the document ---------------
$docEnd = $d.Content.End - 1
$endRange = $d.Range($docEnd, $docEnd)
[void]$endRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>')

Write-Output "OK"
